# FlorianDuruz_JournalDeTravail_CastleDefense.xlsx
# "Correction de la preselection + Documentation sur la preselection
#  (introduction, enjeux, diagramme)"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 44 : fill in the missing "Fin" time and the resulting description ---
$ws.Range("C44").Value = 0.48194444444444445

$ws.Range("E44").Value = "Lecture sur les convensions d'écritures UML"

# Hyperlink in H44 pointing to the UML documentation page. Hyperlinks.Add()
# forces a (new) hyperlink cell style onto the target cell, so immediately
# restore the original formatting (style "9", same as every other cell in
# the "Sources" column) by copying it back from a neighbouring cell.
$ws.Hyperlinks.Add($ws.Range("H44"), "https://www.uml-diagrams.org/abstraction.html")
$ws.Range("H42").Copy()
$ws.Range("H44").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("F44").Value = "Lecture sur la convention concernant les abstract et les interfaces et leur représentation dans un schéma"

# --- Row 45 : new journal entry ---
$ws.Range("A45").Value = 44692
$ws.Range("B45").Value = 0.48194444444444445
$ws.Range("C45").Value = 0.51041666666666663
$ws.Range("E45").Value = "Diagramme de classe" + [char]10 + "Régiment"

# --- Row 46 : new journal entry ---
$ws.Range("A46").Value = 44692
$ws.Range("B46").Value = 0.5625
$ws.Range("C46").Value = 0.56944444444444442
$ws.Range("E46").Value = "BUG avec certains régiments"
$ws.Range("F46").Value = "L'interface Repassait dans tout le registre, ce qui réinitialisait les valeur des autres régiments, précédemment enregistré"

# --- Row 47 : new journal entry ---
$ws.Range("A47").Value = 44692
$ws.Range("B47").Value = 0.56944444444444442
$ws.Range("C47").Value = 0.6333333333333333
$ws.Range("E47").Value = "Documentation:" + [char]10 + "Preselection" + [char]10 + "Diagram de flux"

# --- Row 48 : new journal entry ---
$ws.Range("A48").Value = 44692
$ws.Range("B48").Value = 0.6333333333333333
$ws.Range("C48").Value = 0.70486111111111116
$ws.Range("E48").Value = "Documentation Selection"

# --- Update the saved view position/selection ---
$excel.ActiveWindow.ScrollRow = 43
$ws.Range("C49").Select()
